$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style captured once from a cell known to use the default (unstyled) format,
# so any Text-format detour on the Price column can be reverted cleanly.
$normalStyle = $ws.Range("D9").Style

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.859.67'
$ws.Range("D2").Style = $normalStyle
$ws.Range("E2").Value = '  -3.15%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.618.63'
$ws.Range("D3").Style = $normalStyle
$ws.Range("E3").Value = '  -3.27%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("D4").Style = $normalStyle
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.10'
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = '  -1.72%  '

# Row 6
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = '  +0.10%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3935'
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3846'
$ws.Range("D8").Style = $normalStyle
$ws.Range("E8").Value = '  -2.52%  '

# Row 9
$ws.Range("E9").Value = '  +0.17%  '

# Row 10
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.64'
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = '  -2.30%  '

# Row 11
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.363'
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = '  -2.59%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08448'
$ws.Range("D12").Style = $normalStyle
$ws.Range("E12").Value = '  -2.23%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.83'
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = '  -5.42%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.046'
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = '  -3.80%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.568'
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = '  -1.69%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001279'
$ws.Range("D16").Style = $normalStyle
$ws.Range("E16").Value = '  -2.81%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.596.74'
$ws.Range("D17").Style = $normalStyle
$ws.Range("E17").Value = '  -4.99%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.89'
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = '  -0.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06932'
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = '  -1.15%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.07'
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = '  -4.77%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.812'
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = '  -3.77%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.39'
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = '  -3.86%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.859.04'
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = '  -3.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.451'
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = '  +4.56%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.856'
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = '  +2.51%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.24'
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = '  -3.29%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '156.51'
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = '  -2.23%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '139.82'
$ws.Range("D29").Style = $normalStyle
$ws.Range("E29").Value = '  -4.20%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.272'
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = '  -9.81%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.846'
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = '  -5.82%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.509'
$ws.Range("D32").Style = $normalStyle
$ws.Range("E32").Value = '  +0.67%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.793.54'
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = '  -3.49%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08130'
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = '  -1.55%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9778'
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = '  -1.40%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02887'
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = '  -6.32%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.582'
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = '  -5.36%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2669'
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = '  -5.10%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09145'
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = '  -4.82%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.33'
$ws.Range("D40").Style = $normalStyle
$ws.Range("E40").Value = '  +0.38%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.56'
$ws.Range("D41").Style = $normalStyle
$ws.Range("E41").Value = '  +0.62%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.426'
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = '  -6.26%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7511'
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = '  -4.93%  '

# Row 44
$ws.Range("E44").Value = '  -3.52%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6895'
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = '  -2.79%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.470'
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = '  -3.58%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.069'
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = '  -2.46%  '

# Row 48
$ws.Range("E48").Value = '  -0.04%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08233'
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = '  -5.03%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.94'
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = '  -2.92%  '

# Row 51
$ws.Range("E51").Value = '  -8.65%  '
